$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value (45180 = 2023-09-11) for
# every data row (rows 2 through 525). The update moves that date forward
# by one day (45181 = 2023-09-12) for all of them.
$range = $ws.Range("C2:C525")
$range.Value = 45181
